$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the value column as Text for every cell that would otherwise be
# auto-detected as a number/date by Excel's input parser (dates, numeric
# strings such as house/postal numbers) so the literal text is preserved.
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B4:B9").NumberFormat = "@"

# --- Header row -----------------------------------------------------
$ws.Range("A1").Value = "Daten"
$ws.Range("B1").Value = "Wert"

# --- Personalnummer / Neuer Eintrag rows (green band) ----------------
$ws.Range("A2").Value = "Personalnummer"
$ws.Range("B2").Value = "M100001"
$ws.Range("A3").Value = "Neuer Eintrag gültig ab:"
$ws.Range("B3").Value = "12.12.1992"

# --- Address rows (orange band) ---------------------------------------
$ws.Range("A4").Value = "Strasse"
$ws.Range("B4").Value = "neue Straße"
$ws.Range("A5").Value = "Hausnummer"
$ws.Range("B5").Value = "42"
$ws.Range("A6").Value = "Postleitzahl"
$ws.Range("B6").Value = "10369"
$ws.Range("A7").Value = "Stadt"
$ws.Range("B7").Value = "Berlin"
$ws.Range("A8").Value = "Region"
$ws.Range("B8").Value = "Berlin"
$ws.Range("A9").Value = "Land"
$ws.Range("B9").Value = "Deutschland"

# --- Formatting: bold header, coloured bands --------------------------
$ws.Range("A1:B1").Font.Bold = $true

$ws.Range("A2:B3").Interior.ThemeColor = 10
$ws.Range("A4:A9").Interior.ThemeColor = 6
$ws.Range("B4:B9").Interior.ThemeColor = 6

# --- Column widths (auto-fit to content) --------------------------------
$ws.Columns.Item(1).ColumnWidth = 55.11
$ws.Columns.Item(2).ColumnWidth = 25.67

# --- Page setup ------------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
